# Auto-generated Excel COM-interop script to apply scheduled runner updates
# to the Sheets workbook, based on the authoritative XML diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1964.8387
$ws.Range("J17").Value = 1964.8387
$ws.Range("L17").Value = 5894.5161
$ws.Range("N17").Value = -6230.5161

$ws.Range("H42").Value = 804.8889
$ws.Range("I42").Value = 40.666668
$ws.Range("K42").Value = 122.000004
$ws.Range("M42").Value = 107.999996

$ws.Range("H82").Value = 2853
$ws.Range("I82").Value = 2853
$ws.Range("K82").Value = 8559
$ws.Range("M82").Value = -8153

$ws.Range("H85").Value = 2853
$ws.Range("I85").Value = 2853
$ws.Range("K85").Value = 8559
$ws.Range("M85").Value = -7155

$ws.Range("H98").Value = 2131.72
$ws.Range("I98").Value = 2111.8235
$ws.Range("K98").Value = 2111.8235
$ws.Range("M98").Value = -613.8235

$ws.Range("H101").Value = 2143.2307
$ws.Range("J101").Value = 4979.6665
$ws.Range("L101").Value = 14938.9995
$ws.Range("N101").Value = -18182.9995

$ws.Range("H122").Value = 2131.72
$ws.Range("I122").Value = 2111.8235
$ws.Range("K122").Value = 6335.470499999999
$ws.Range("M122").Value = -3885.470499999999

$ws.Range("H129").Value = 1784
$ws.Range("I129").Value = 1386.4286
$ws.Range("J129").Value = 3175.5
$ws.Range("K129").Value = 4159.2858
$ws.Range("L129").Value = 9526.5
$ws.Range("M129").Value = 840.7142000000003
$ws.Range("N129").Value = -19526.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11367621
$ws.Range("I32").Value = 11631868
$ws.Range("K32").Value = 11631868
$ws.Range("M32").Value = -11631581

$ws.Range("H45").Value = 3915.25
$ws.Range("I45").Value = 3869.125
$ws.Range("K45").Value = 3869.125
$ws.Range("M45").Value = -3492.125

$ws.Range("H132").Value = 7625
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 7625
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 22875
$ws.Range("N132").Value = -27935
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 76374
$ws.Range("I99").Value = 104126.8
$ws.Range("J99").Value = 6992
$ws.Range("K99").Value = 104126.8
$ws.Range("L99").Value = 6992
$ws.Range("M99").Value = -102628.8
$ws.Range("N99").Value = -9988

$ws.Range("H100").Value = 9969.5
$ws.Range("J100").Value = 9969.5
$ws.Range("L100").Value = 9969.5
$ws.Range("N100").Value = -12133.5

$ws.Range("H103").Value = 4999
$ws.Range("J103").Value = 4999
$ws.Range("L103").Value = 4999
$ws.Range("N103").Value = -7343

$ws.Range("H105").Value = 2665.8333
$ws.Range("I105").Value = 2665.8333
$ws.Range("K105").Value = 2665.8333
$ws.Range("M105").Value = -918.8332999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5236.0337
$ws.Range("I31").Value = 2211.5715
$ws.Range("J31").Value = 7967.8066
$ws.Range("K31").Value = 2211.5715
$ws.Range("L31").Value = 7967.8066
$ws.Range("M31").Value = -1916.5715
$ws.Range("N31").Value = -8557.8066

$ws.Range("H34").Value = 5236.0337
$ws.Range("I34").Value = 2211.5715
$ws.Range("J34").Value = 7967.8066
$ws.Range("K34").Value = 2211.5715
$ws.Range("L34").Value = 7967.8066
$ws.Range("M34").Value = -2009.5715
$ws.Range("N34").Value = -8371.8066

$ws.Range("H122").Value = 3026.8096
$ws.Range("I122").Value = 1588.6154
$ws.Range("K122").Value = 4765.8462
$ws.Range("M122").Value = -2315.8462

$ws.Range("H132").Value = 6589532
$ws.Range("I132").Value = 11389.147
$ws.Range("J132").Value = 62503748
$ws.Range("K132").Value = 34167.44100000001
$ws.Range("L132").Value = 187511244
$ws.Range("M132").Value = -31637.44100000001
$ws.Range("N132").Value = -187516304

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6928.769
$ws.Range("I3").Value = 1512.5555
$ws.Range("K3").Value = 4537.666499999999
$ws.Range("M3").Value = -4425.666499999999

$ws.Range("H92").Value = 1073.3
$ws.Range("I92").Value = 912.6429000000001
$ws.Range("K92").Value = 2737.9287
$ws.Range("M92").Value = -1489.9287

$ws.Range("H114").Value = 2001.6666
$ws.Range("I114").Value = 146.61539
$ws.Range("J114").Value = 6824.8
$ws.Range("K114").Value = 439.84617
$ws.Range("L114").Value = 20474.4
$ws.Range("M114").Value = 2814.15383
$ws.Range("N114").Value = -26982.4

$ws.Range("H131").Value = 14050.111
$ws.Range("J131").Value = 17821.857
$ws.Range("L131").Value = 53465.571
$ws.Range("N131").Value = -63545.571

$ws.Range("H133").Value = 3666.2593
$ws.Range("I133").Value = 3285.1904
$ws.Range("K133").Value = 9855.5712
$ws.Range("M133").Value = -4795.5712

$ws.Range("H134").Value = 2107.5
$ws.Range("I134").Value = 2107.5
$ws.Range("K134").Value = 6322.5
$ws.Range("M134").Value = -1252.5

$ws.Range("H138").Value = 2959.111
$ws.Range("I138").Value = 2959.111
$ws.Range("K138").Value = 8877.332999999999
$ws.Range("M138").Value = -3737.332999999999

$ws.Range("H139").Value = 3881.3333
$ws.Range("I139").Value = 2500
$ws.Range("K139").Value = 7500
$ws.Range("M139").Value = -2360

$ws.Range("H140").Value = 4359.95
$ws.Range("I140").Value = 2822.875
$ws.Range("K140").Value = 8468.625
$ws.Range("M140").Value = -3288.625

$ws.Range("H141").Value = 4793
$ws.Range("I141").Value = 4793
$ws.Range("K141").Value = 14379
$ws.Range("M141").Value = -9199

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H102").Value = 2001.6562
$ws.Range("I102").Value = 1309.7693
$ws.Range("K102").Value = 1309.7693
$ws.Range("M102").Value = 312.2307000000001

$ws.Range("H132").Value = 465369.5
$ws.Range("I132").Value = 604036.0600000001
$ws.Range("K132").Value = 1812108.18
$ws.Range("M132").Value = -1809578.18

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 723.8
$ws.Range("J9").Value = 1073
$ws.Range("L9").Value = 1073
$ws.Range("N9").Value = -1521

$ws.Range("H22").Value = 55416.79
$ws.Range("I22").Value = 127318.75
$ws.Range("K22").Value = 127318.75
$ws.Range("M22").Value = -127023.75

$ws.Range("H25").Value = 46599.8
$ws.Range("I25").Value = 47666.332
$ws.Range("J25").Value = 45000
$ws.Range("K25").Value = 47666.332
$ws.Range("L25").Value = 45000
$ws.Range("M25").Value = -47436.332
$ws.Range("N25").Value = -45460

$ws.Range("H27").Value = 55416.79
$ws.Range("I27").Value = 127318.75
$ws.Range("K27").Value = 127318.75
$ws.Range("M27").Value = -127211.75

$ws.Range("H46").Value = 942
$ws.Range("I46").Value = 966
$ws.Range("K46").Value = 966
$ws.Range("M46").Value = -778

$ws.Range("H93").Value = 1970.4286
$ws.Range("I93").Value = 1848.5
$ws.Range("J93").Value = 2081.2727
$ws.Range("K93").Value = 1848.5
$ws.Range("L93").Value = 2081.2727
$ws.Range("M93").Value = -600.5
$ws.Range("N93").Value = -4577.2727

$ws.Range("H100").Value = 8026.5
$ws.Range("I100").Value = 2977.7856
$ws.Range("J100").Value = 16861.75
$ws.Range("K100").Value = 2977.7856
$ws.Range("L100").Value = 16861.75
$ws.Range("M100").Value = -2436.7856
$ws.Range("N100").Value = -17943.75

$ws.Range("H122").Value = 3316.925
$ws.Range("I122").Value = 3063.0344
$ws.Range("K122").Value = 9189.1032
$ws.Range("M122").Value = -6739.1032

$ws.Range("H132").Value = 1237881.1
$ws.Range("I132").Value = 1505964
$ws.Range("J132").Value = 4699.8
$ws.Range("K132").Value = 4517892
$ws.Range("L132").Value = 14099.4
$ws.Range("M132").Value = -4515362
$ws.Range("N132").Value = -19159.4

$ws.Range("H136").Value = 6474.311
$ws.Range("I136").Value = 5631.788
$ws.Range("K136").Value = 16895.364
$ws.Range("M136").Value = -14345.364

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H122").Value = 2992.2
$ws.Range("I122").Value = 2913.5454
$ws.Range("K122").Value = 8740.636200000001
$ws.Range("M122").Value = -6290.636200000001

$ws.Range("H126").Value = 3110.8572
$ws.Range("I126").Value = 2604.6365
$ws.Range("J126").Value = 4967
$ws.Range("K126").Value = 7813.9095
$ws.Range("L126").Value = 14901
$ws.Range("M126").Value = -5343.9095
$ws.Range("N126").Value = -19841

$ws.Range("H132").Value = 4911274
$ws.Range("I132").Value = 5921148.5
$ws.Range("K132").Value = 17763445.5
$ws.Range("M132").Value = -17760915.5
